$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 721, shifting rows 721:743 down to 722:744
$ws.Rows.Item(721).Insert()

# Fill the new row 721 with the new record's data
$ws.Range("A721").Value = 3
$ws.Range("B721").Value = "Femacal de La Calera"
$ws.Range("C721").Value = "Coquimbo"
$ws.Range("D721").Value = 45239
$ws.Range("D721").NumberFormat = $ws.Range("D722").NumberFormat
$ws.Range("E721").Value = 5
$ws.Range("F721").Value = 100112028
$ws.Range("G721").Value = "Sandia"
$ws.Range("H721").Value = "Sin especificar"
$ws.Range("I721").Value = "Primera"
$ws.Range("J721").Value = 160
$ws.Range("K721").Value = 700
$ws.Range("L721").Value = 700
$ws.Range("M721").Value = 700
$ws.Range("N721").Value = "$/kilo (volumen en unidades)"
$ws.Range("O721").Value = "Perú"
$ws.Range("P721").Value = 700
$ws.Range("Q721").Value = 1
$ws.Range("R721").Value = "Hortaliza"
